$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2:D51 and E2:E51 hold text values (e.g. "24.704.12", "  +0.55%  ") in the
# source workbook. Temporarily switch their number format to Text so Excel does
# not reinterpret numeric-looking strings (like "316.19" or "1.003") as real
# numbers when the new values are assigned below.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

# Updated coin price (column D) and 1h volume change (column E) figures.
$ws.Cells.Item(2, 4).Value = "24.704.12"
$ws.Cells.Item(2, 5).Value = "  +0.55%  "
$ws.Cells.Item(3, 4).Value = "1.699.40"
$ws.Cells.Item(3, 5).Value = "  +0.32%  "
$ws.Cells.Item(4, 5).Value = "  +0.21%  "
$ws.Cells.Item(5, 4).Value = "316.19"
$ws.Cells.Item(5, 5).Value = "  +0.09%  "
$ws.Cells.Item(6, 4).Value = "1.003"
$ws.Cells.Item(6, 5).Value = "  +0.21%  "
$ws.Cells.Item(7, 4).Value = "0.3929"
$ws.Cells.Item(7, 5).Value = "  -0.18%  "
$ws.Cells.Item(8, 4).Value = "0.4052"
$ws.Cells.Item(8, 5).Value = "  +1.04%  "
$ws.Cells.Item(9, 4).Value = "1.519"
$ws.Cells.Item(9, 5).Value = "  -0.34%  "
$ws.Cells.Item(10, 4).Value = "1.005"
$ws.Cells.Item(10, 5).Value = "  +0.31%  "
$ws.Cells.Item(11, 4).Value = "52.89"
$ws.Cells.Item(11, 5).Value = "  +0.15%  "
$ws.Cells.Item(12, 4).Value = "0.08877"
$ws.Cells.Item(12, 5).Value = "  +1.63%  "
$ws.Cells.Item(13, 4).Value = "7.421"
$ws.Cells.Item(13, 5).Value = "  +3.09%  "
$ws.Cells.Item(14, 4).Value = "23.69"
$ws.Cells.Item(14, 5).Value = "  +2.46%  "
$ws.Cells.Item(15, 4).Value = "8.115"
$ws.Cells.Item(15, 5).Value = "  +7.37%  "
$ws.Cells.Item(16, 4).Value = "0.00001321"
$ws.Cells.Item(16, 5).Value = "  +0.33%  "
$ws.Cells.Item(17, 4).Value = "1.703.03"
$ws.Cells.Item(17, 5).Value = "  +0.55%  "
$ws.Cells.Item(18, 4).Value = "99.46"
$ws.Cells.Item(18, 5).Value = "  -0.23%  "
$ws.Cells.Item(19, 4).Value = "0.07063"
$ws.Cells.Item(19, 5).Value = "  +0.21%  "
$ws.Cells.Item(20, 4).Value = "19.78"
$ws.Cells.Item(20, 5).Value = "  +0.58%  "
$ws.Cells.Item(21, 4).Value = "7.065"
$ws.Cells.Item(21, 5).Value = "  +2.95%  "
$ws.Cells.Item(22, 4).Value = "1.006"
$ws.Cells.Item(22, 5).Value = "  +0.55%  "
$ws.Cells.Item(23, 4).Value = "14.75"
$ws.Cells.Item(23, 5).Value = "  +5.13%  "
$ws.Cells.Item(24, 4).Value = "24.696.17"
$ws.Cells.Item(24, 5).Value = "  +0.54%  "
$ws.Cells.Item(25, 4).Value = "3.139"
$ws.Cells.Item(25, 5).Value = "  +4.49%  "
$ws.Cells.Item(26, 4).Value = "2.351"
$ws.Cells.Item(26, 5).Value = "  +1.29%  "
$ws.Cells.Item(27, 4).Value = "22.65"
$ws.Cells.Item(27, 5).Value = "  +1.41%  "
$ws.Cells.Item(28, 4).Value = "164.69"
$ws.Cells.Item(28, 5).Value = "  +2.82%  "
$ws.Cells.Item(29, 4).Value = "8.840"
$ws.Cells.Item(29, 5).Value = "  +19.06%  "
$ws.Cells.Item(30, 4).Value = "135.74"
$ws.Cells.Item(30, 5).Value = "  +0.98%  "
$ws.Cells.Item(31, 4).Value = "5.153"
$ws.Cells.Item(31, 5).Value = "  -1.26%  "
$ws.Cells.Item(32, 4).Value = "0.09017"
$ws.Cells.Item(32, 5).Value = "  +6.02%  "
$ws.Cells.Item(33, 4).Value = "7.664"
$ws.Cells.Item(33, 5).Value = "  +6.30%  "
$ws.Cells.Item(34, 4).Value = "1.069"
$ws.Cells.Item(34, 5).Value = "  -2.51%  "
$ws.Cells.Item(35, 4).Value = "0.03011"
$ws.Cells.Item(35, 5).Value = "  +10.00%  "
$ws.Cells.Item(36, 4).Value = "1.961"
$ws.Cells.Item(36, 5).Value = "  +0.17%  "
$ws.Cells.Item(37, 4).Value = "0.2754"
$ws.Cells.Item(37, 5).Value = "  +1.17%  "
$ws.Cells.Item(38, 4).Value = "10.97"
$ws.Cells.Item(38, 5).Value = "  -4.22%  "
$ws.Cells.Item(39, 4).Value = "14.43"
$ws.Cells.Item(39, 5).Value = "  +0.07%  "
$ws.Cells.Item(40, 4).Value = "0.09218"
$ws.Cells.Item(40, 5).Value = "  +1.91%  "
$ws.Cells.Item(41, 4).Value = "1.467"
$ws.Cells.Item(41, 5).Value = "  +0.28%  "
$ws.Cells.Item(42, 4).Value = "0.7684"
$ws.Cells.Item(42, 5).Value = "  +0.04%  "
$ws.Cells.Item(43, 4).Value = "16.02"
$ws.Cells.Item(43, 5).Value = "  +3.99%  "
$ws.Cells.Item(44, 4).Value = "0.7174"
$ws.Cells.Item(44, 5).Value = "  -0.19%  "
$ws.Cells.Item(45, 4).Value = "2.584"
$ws.Cells.Item(45, 5).Value = "  +1.86%  "
$ws.Cells.Item(46, 4).Value = "4.220"
$ws.Cells.Item(46, 5).Value = "  +0.30%  "
$ws.Cells.Item(47, 5).Value = "  +0.12%  "
$ws.Cells.Item(48, 4).Value = "1.352"
$ws.Cells.Item(48, 5).Value = "  +2.29%  "
$ws.Cells.Item(49, 4).Value = "139.89"
$ws.Cells.Item(49, 5).Value = "  -0.66%  "
$ws.Cells.Item(50, 4).Value = "0.07972"
$ws.Cells.Item(50, 5).Value = "  -0.38%  "
$ws.Cells.Item(51, 4).Value = "89.78"
$ws.Cells.Item(51, 5).Value = "  +1.74%  "

# Restore the original (default/"Normal") cell style now that the text values
# are in place, so no lasting formatting change is introduced.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
